# Apply updated loading_percent values for the "380 kV" case
# The workbook has a single worksheet with a results table:
#   Row 1: header (column indices 0-13)
#   Rows 2-25: data rows, column A = case index 0..23
#   Data columns used: B, C, D, E, F, I, M, N, O
# This script overwrites the numeric values in those columns for rows 2-25
# with the updated simulation results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 11.17622703276411, 8.448252368653652, 3.852438476026334, 11.87252651299298, 20.85295102378896, 17.4299331331869, 14.15154519726645, 16.33347107265294, 18.4545401647835),
    @(3, 10.61141056364529, 8.010231711293601, 3.805593125082941, 11.75994438436754, 20.77413510526154, 17.51116580657229, 13.87218533177987, 16.3846106211999, 18.46272471924559),
    @(4, 10.24956465618191, 7.727085927984708, 3.776170009103386, 11.69489938987281, 20.73284152840559, 17.56551027831302, 13.7010171826731, 16.41776382012989, 18.47334203089321),
    @(5, 10.09849470632513, 7.608197006260046, 3.764020639967602, 11.66944737282207, 20.71781252560319, 17.58877572524529, 13.63145501701311, 16.4317158433105, 18.47907199672589),
    @(6, 10.07319683661219, 7.588246105740095, 3.761993829930942, 11.66528553646765, 20.71542592454843, 17.59270644972904, 13.61991870714797, 16.43405927928135, 18.48010810958365),
    @(7, 10.24754166594828, 7.725496644043865, 3.776006794118875, 11.69455183079318, 20.73263154478516, 17.56581951583464, 13.70007813978195, 16.41795019180918, 18.47341362966697),
    @(8, 10.98469371046189, 8.300205396035516, 3.836427139027135, 11.83287889940988, 20.82431050332281, 17.45701275459161, 14.0552065066052, 16.35074065412261, 18.4562001744257),
    @(9, 12.30480177999834, 9.312272240194075, 3.949389528930546, 12.13515722279192, 21.05971489511857, 17.27925375878475, 14.75018965945444, 16.23281207741388, 18.46689838648294),
    @(10, 13.19169841628284, 9.983448307774273, 4.028648332275946, 12.37417730824791, 21.26543920845329, 17.17058806886879, 15.25440062742452, 16.15456584690867, 18.50191149302953),
    @(11, 13.57615563574226, 10.27273903397859, 4.063820199713073, 12.48613609818148, 21.3658494193451, 17.12596383555264, 15.48129520752819, 16.12078046575481, 18.52372566244047),
    @(12, 13.71894466047142, 10.37996137389597, 4.077005891149653, 12.52895625092656, 21.40482522645563, 17.10976119223357, 15.56677328528756, 16.10824611085534, 18.53282991719347),
    @(13, 13.68831777056039, 10.35697283830784, 4.074172132284549, 12.51971596567684, 21.39638920126547, 17.11321971927105, 15.54838502997231, 16.11093408131475, 18.53083168487972),
    @(14, 13.58795935021584, 10.28160699237788, 4.064907696752655, 12.489650677855, 21.36903706267591, 17.12461687010644, 15.48833687491144, 16.11974406144155, 18.5244577856629),
    @(15, 13.52612100863865, 10.23513982666674, 4.059215447736046, 12.47128880263167, 21.35240625849971, 17.13168866248476, 15.45149561133064, 16.12517419363795, 18.52066336239788),
    @(16, 13.1661851364815, 9.96421802120752, 4.026331465215094, 12.36692205285478, 21.25901209109971, 17.17360147919109, 15.23951527139821, 16.15681014439134, 18.50060412554835),
    @(17, 12.94046039803229, 9.793894094314975, 4.005927538100516, 12.3036946136285, 21.20344607562684, 17.20054812955158, 15.10877572963401, 16.17668063482697, 18.48980448738369),
    @(18, 12.80884445391961, 9.694420261452599, 3.99410894317505, 12.26763410885655, 21.17213025961129, 17.21649946677227, 15.03334984215856, 16.18827994983412, 18.48414712445017),
    @(19, 12.76397716079555, 9.660481841790654, 3.990093339695833, 12.25547845482136, 21.16163876186693, 17.22197787768609, 15.00777548948644, 16.19223655629738, 18.48232690742199),
    @(20, 12.96467443719919, 9.812181606561991, 4.008108184994612, 12.31039387674395, 21.20929469352028, 17.19763276726133, 15.1227174108028, 16.17454776355858, 18.4908967756029),
    @(21, 13.6175134053538, 10.30380702810215, 4.067632546655158, 12.49847040478033, 21.37704543452516, 17.12125033972348, 15.50598713183452, 16.11714932139379, 18.52630708402713),
    @(22, 14.02785902395911, 10.61155010428555, 4.105756505171384, 12.62384034453094, 21.49221814268567, 17.07538631251686, 15.75386339120493, 16.08114799037169, 18.55436545671321),
    @(23, 13.81036167651891, 10.44854833299779, 4.08548223238733, 12.55671725300026, 21.43025163653811, 17.09949229214111, 15.62183344261978, 16.10022447652569, 18.53894158636472),
    @(24, 12.95373300801968, 9.803918655917322, 4.007122588490117, 12.30736423416212, 21.20664856931424, 17.19894937198979, 15.11641518613719, 16.17551148789399, 18.49040123400433),
    @(25, 11.96185823515192, 9.051059709466712, 3.919460150961907, 12.05025863478781, 20.9901840011887, 17.32350665063151, 14.56291480160232, 16.263236294325, 18.45923558686937)
)

foreach ($row in $data) {
    $r = $row[0]
    # Columns B,C,D,E,F map to indices 2..6
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    # Column I maps to index 9
    $ws.Cells.Item($r, 9).Value = $row[6]
    # Columns M,N,O map to indices 13..15
    $ws.Cells.Item($r, 13).Value = $row[7]
    $ws.Cells.Item($r, 14).Value = $row[8]
    $ws.Cells.Item($r, 15).Value = $row[9]
}

